# fix: fixed formatting when scrapping floating point numbers
#
# A prior bulk "clean up the scraped numbers" pass on the "Importe"
# column (H) used a naive comma/period swap that also clobbered a
# handful of "Razon social" (name) cells in column E/F which happened
# to contain a literal comma (used there as a list-of-people separator)
# or a "S.H." abbreviation. This restores the intended, uniform rule
# ("drop thousands-separator periods, turn the decimal comma into a
# decimal period") only where it belongs, and reverts the name cells
# back to using a period as their separator instead of a comma.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Razon social / name cells that used a comma as a list separator.
#    These are plain text - just rewrite them directly.
# ---------------------------------------------------------------------
$ws.Range("E40").Value  = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E105").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"

$ws.Range("E55").Value  = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Range("E199").Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"

$ws.Range("E189").Value = "RICCOTTI. MARIANA EDITH"

$ws.Range("E217").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"

$ws.Range("E290").Value = "ALVAREZ. RENZO JOEL"
$ws.Range("F290").Value = "ALVAREZ. RENZO JOEL"

# ---------------------------------------------------------------------
# 2) "Importe" column (H2:H298): these amounts are stored as literal
#    text (e.g. "1.056.600,00"), not real numbers. Simply assigning the
#    reformatted text back (e.g. "1056600.00") would make the engine
#    auto-detect a real number and silently drop the trailing ".00".
#    So: compute the fixed text in a scratch column via SUBSTITUTE(),
#    copy it, and paste-special VALUES ONLY back onto the Importe
#    column - that carries over the literal text without retyping it
#    as a number and without touching any cell's number format/style.
# ---------------------------------------------------------------------
$target = $ws.Range("H2:H298")
$scratch = $ws.Range("Z2:Z298")

$scratch.Formula = '=SUBSTITUTE(SUBSTITUTE(H2,".",""),",",".")'

$scratch.Copy()
$target.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0

$scratch.ClearContents()
